# Add 'Place Bid' to Test Procedure column (column E) for all 'add bid' testcases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 19; $r -le 35; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 + ". Place Bid."
}

$ws.Range("E36").Select()
